$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, bordered, centered) from an existing header cell
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-16
$IValues = @(1,1,1,1,1,1,7,1,1,1,1,6,7,5,1)
$JValues = @(4,5,4,2,6,6,9,2,1,2,3,7,8,6,2)

for ($r = 0; $r -lt 15; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $IValues[$r]
    $ws.Cells.Item($row, 10).Value = $JValues[$r]
}
